# "Continue test on dataset 2a"
# Fills in the previously-empty Fold 5 / Fold 6 (columns F / G) accuracy
# results for the "D2A (Mixed Up)(128)" sheet, and updates the view state
# (scroll position / selection) on that sheet and on "D2A (Mixed Up)(64)".
# All AVERAGE / MAX formulas (cols B, C) and the derived *100 tables
# (rows 34-57) recalc automatically from these inputs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "D2A (Mixed Up)(64)" - just a scroll-position change
# ---------------------------------------------------------------------
$wsMixed64 = $wb.Worksheets.Item("D2A (Mixed Up)(64)")
$wsMixed64.Activate()
$excel.ActiveWindow.ScrollRow = 19
$wsMixed64.Range("C34").Select()

# ---------------------------------------------------------------------
# Sheet "D2A (Mixed Up)(128)" - new fold data + view state
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("D2A (Mixed Up)(128)")
$ws.Activate()

# --- first block (rows 6-14) ---
$ws.Range("F6").Value  = 0.84722200000000003
$ws.Range("G6").Value  = 0.85763900000000004

$ws.Range("F7").Value  = 0.77083299999999999
$ws.Range("G7").Value  = 0.72916700000000001

$ws.Range("F8").Value  = 0.89583299999999999
$ws.Range("G8").Value  = 0.89930600000000005

$ws.Range("F9").Value  = 0.734375
$ws.Range("G9").Value  = 0.72395799999999999

$ws.Range("F10").Value = 0.84722200000000003
$ws.Range("G10").Value = 0.80902799999999997

$ws.Range("F11").Value = 0.86805600000000005
$ws.Range("G11").Value = 0.83680600000000005

$ws.Range("F12").Value = 0.86458299999999999
$ws.Range("G12").Value = 0.86805600000000005

$ws.Range("F13").Value = 0.86805600000000005
$ws.Range("G13").Value = 0.87152799999999997

$ws.Range("F14").Value = 0.92708299999999999
$ws.Range("G14").Value = 0.90277799999999997

# --- second block (rows 19-27) ---
$ws.Range("F19").Value = 0.84375
$ws.Range("G19").Value = 0.82638900000000004

$ws.Range("F20").Value = 0.77083299999999999
$ws.Range("G20").Value = 0.71180600000000005

$ws.Range("F21").Value = 0.88541700000000001
$ws.Range("G21").Value = 0.875

$ws.Range("F22").Value = 0.77083299999999999
$ws.Range("G22").Value = 0.72916700000000001

$ws.Range("F23").Value = 0.83333299999999999
$ws.Range("G23").Value = 0.81597200000000003

$ws.Range("F24").Value = 0.84722200000000003
$ws.Range("G24").Value = 0.84027799999999997

$ws.Range("F25").Value = 0.86805600000000005
$ws.Range("G25").Value = 0.875

$ws.Range("F26").Value = 0.87152799999999997
$ws.Range("G26").Value = 0.89236099999999996

$ws.Range("F27").Value = 0.91319399999999995
$ws.Range("G27").Value = 0.90277799999999997

# --- view state: scrolled to row 7, selection at G62 ---
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("G62").Select()
